# Add new columns I (I0) and J (IF) to the worksheet, mirroring the
# existing header style used for column H (header row style index "1",
# i.e. the style currently applied to cells B1:H1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style from H1 (bold, bordered, centered) onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data rows (2-29): column I then column J ---
$data = @{
    2  = @(6, 7)
    3  = @(4, 5)
    4  = @(5, 8)
    5  = @(1, 5)
    6  = @(1, 7)
    7  = @(1, 6)
    8  = @(1, 5)
    9  = @(1, 4)
    10 = @(1, 5)
    11 = @(1, 5)
    12 = @(1, 6)
    13 = @(1, 6)
    14 = @(1, 5)
    15 = @(1, 7)
    16 = @(1, 7)
    17 = @(1, 3)
    18 = @(1, 7)
    19 = @(1, 5)
    20 = @(1, 5)
    21 = @(1, 4)
    22 = @(1, 2)
    23 = @(1, 6)
    24 = @(1, 6)
    25 = @(4, 8)
    26 = @(1, 4)
    27 = @(1, 3)
    28 = @(1, 2)
    29 = @(1, 1)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]   # column I
    $ws.Cells.Item($row, 10).Value = $vals[1]  # column J
}
